$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove stale hyperlinks up front (row deletion / re-targeting below) ---
$ws.Hyperlinks.Delete()

# --- Update data rows 2-13 with the freshly scraped listings ---
$ws.Range('A2').Value = '2025-09-26 06:29:13'
$ws.Range('B2').Value = 'LINExChatGPTx美容室向け予約Bot (仕様書、契約書あり)'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('F2').Value = 'https://www.lancers.jp/work/detail/5400801'
$ws.Range('G2').Value = 445
$ws.Range('H2').Value = '🔥GPT,ChatGPT ★bot'

$ws.Range('A3').Value = '2025-09-26 06:29:13'
$ws.Range('B3').Value = '自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('F3').Value = 'https://www.lancers.jp/work/detail/5389460'
$ws.Range('G3').Value = 305
$ws.Range('H3').Value = '🔥Python ◆開発 ○PHP'

$ws.Range('A4').Value = '2025-09-26 06:29:13'
$ws.Range('B4').Value = 'システムの開発補助や運営サポート【フルリモート×長期】'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('F4').Value = 'https://www.lancers.jp/work/detail/5323359'
$ws.Range('G4').Value = 83
$ws.Range('H4').Value = '◆開発'

$ws.Range('A5').Value = '2025-09-26 06:29:13'
$ws.Range('B5').Value = '【急募】LLMによるMCP(Model Context Protocol)でのExcel操作機能開発'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('F5').Value = 'https://www.lancers.jp/work/detail/5400689'
$ws.Range('G5').Value = 75
$ws.Range('H5').Value = '◆開発'

$ws.Range('A6').Value = '2025-09-26 06:29:13'
$ws.Range('B6').Value = '【急募】音源ライセンス販売サイトのMVP構築依頼'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('F6').Value = 'https://www.lancers.jp/work/detail/5400763'
$ws.Range('G6').Value = 45
$ws.Range('H6').Value = '◇サイト'

$ws.Range('A7').Value = '2025-09-26 06:29:13'
$ws.Range('B7').Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('F7').Value = 'https://www.lancers.jp/work/detail/5016989'
$ws.Range('G7').Value = 33
$ws.Range('H7').Value = '○WordPress'

$ws.Range('A8').Value = '2025-09-26 06:29:13'
$ws.Range('B8').Value = 'eBayテラピークでのキーワード検索結果等の取得するためのシステム制作'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('F8').Value = 'https://www.lancers.jp/work/detail/5390238'
$ws.Range('G8').Value = 33
$ws.Range('H8').ClearContents()

$ws.Range('A9').Value = '2025-09-26 06:29:13'
$ws.Range('B9').Value = 'Drupal関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5400683'
$ws.Range('G9').Value = 25
$ws.Range('H9').ClearContents()

$ws.Range('A10').Value = '2025-09-26 06:29:13'
$ws.Range('B10').Value = '金融関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)'
$ws.Range('C10').Value = 'システム開発'
$ws.Range('D10').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E10').Value = '期限情報なし'
$ws.Range('F10').Value = 'https://www.lancers.jp/work/detail/5400681'
$ws.Range('G10').Value = 25
$ws.Range('H10').ClearContents()

$ws.Range('A11').Value = '2025-09-26 06:29:13'
$ws.Range('B11').Value = '【急募】東京でのWeb制作プロジェクトに参加しませんか?'
$ws.Range('C11').Value = 'システム開発'
$ws.Range('D11').Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range('E11').Value = '期限情報なし'
$ws.Range('F11').Value = 'https://www.lancers.jp/work/detail/5400965'
$ws.Range('G11').Value = 18
$ws.Range('H11').ClearContents()

$ws.Range('A12').Value = '2025-09-26 06:29:13'
$ws.Range('B12').Value = '限定公開 PR 限定公開の仕事'
$ws.Range('C12').Value = 'システム開発'
$ws.Range('D12').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E12').Value = '期限情報なし'
$ws.Range('F12').Value = 'https://www.lancers.jp/work/detail/5399347'
$ws.Range('G12').Value = 13
$ws.Range('H12').ClearContents()

$ws.Range('A13').Value = '2025-09-26 06:29:13'
$ws.Range('B13').Value = '【急募】スーパードルフィーの洋服をオーダーメイドで作成希望'
$ws.Range('C13').Value = 'システム開発'
$ws.Range('D13').Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range('E13').Value = '期限情報なし'
$ws.Range('F13').Value = 'https://www.lancers.jp/work/detail/5400988'
$ws.Range('G13').Value = 10
$ws.Range('H13').ClearContents()

# --- Drop the now-stale trailing rows (old rows 14-20) ---
$ws.Range('A14:H20').EntireRow.Delete()

# --- Re-create hyperlinks for F2:F13 against their new target URLs ---
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5400801') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5389460') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5323359') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5400689') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5400763') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5016989') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5390238') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5400683') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F10'), 'https://www.lancers.jp/work/detail/5400681') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F11'), 'https://www.lancers.jp/work/detail/5400965') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F12'), 'https://www.lancers.jp/work/detail/5399347') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F13'), 'https://www.lancers.jp/work/detail/5400988') | Out-Null

# --- Column width tweaks (attribute width = ColumnWidth + 5/6) ---
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6
$ws.Columns.Item(8).ColumnWidth = 19 - 5/6
